# "fixes for api & web demo tests"
# Replace the three demo phone rows on the GSMArena sheet with the
# Galaxy S10+ / Galaxy S10 / Galaxy View2 data, give the new D2 model
# cell its own (Menlo / dark-grey) font, and leave the GSMArena sheet
# selected/active (matching the refreshed view state in the workbook).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # GSMArena
$ws2 = $wb.Worksheets.Item(2)   # Calculator

# --- Row 2: Galaxy S10+ ---------------------------------------------------
$ws1.Range("D2").Value = "Galaxy S10+"
$ws1.Range("E2").Value = "6.4"""
$ws1.Range("F2").Value = "16MP"
$ws1.Range("G2").Value = "12GB RAM"
$ws1.Range("H2").Value = "4100mAh"

# --- Row 3: Galaxy S10 -----------------------------------------------------
$ws1.Range("D3").Value = "Galaxy S10"
$ws1.Range("E3").Value = "6.1"""
$ws1.Range("F3").Value = "16MP"
$ws1.Range("G3").Value = "8GB RAM"
$ws1.Range("H3").Value = "3400mAh"

# --- Row 4: Galaxy View2 ----------------------------------------------------
$ws1.Range("D4").Value = "Galaxy View2"
$ws1.Range("E4").Value = "17.3"""
$ws1.Range("F4").Value = "NO"
$ws1.Range("G4").Value = "3GB RAM"
$ws1.Range("H4").Value = "12000mAh"

# New cell-level font (Menlo, #222222) applied to the new model name in D2
$d2 = $ws1.Range("D2")
$d2.ClearFormats()
$d2.Font.Color = 2236962   # RGB(34,34,34) = 0x222222
$d2.Font.Name  = "Menlo"
$d2.Value = "Galaxy S10+"  # ClearFormats() resets the value, restore it

# GSMArena becomes the selected/active sheet, cursor parked at H5
$ws1.Activate() | Out-Null
$ws1.Range("H5").Select() | Out-Null
